# Table 3 note: rephrase the "Cluster location" sentence.
#   before: "... most significant electrode; clusters often extend into additional regions. "
#   after:  "... most significant electrode per cluster. "
# The "_GoBack" bookmark (currently wrapping the whole sentence) ends up
# collapsed/empty, sitting right after "cluster" and before the new ". ".

$d = $word.ActiveDocument

# --- 1) Insert " per" right after "... most significant electrode" ---------
$perRng = $d.Content
$perRng.Find.Execute("most significant electrode", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$perRng.Collapse(0)
$perRng.InsertAfter(" per")

# --- 2) Insert " cluster" right after the newly-added " per" ---------------
$clusterRng = $perRng.Duplicate
$clusterRng.Collapse(0)
$clusterRng.InsertAfter(" cluster")

# Nudge formatting off/on (no visible change) so each inserted chunk keeps
# its own run instead of silently re-merging with its neighbour; do the
# right-most range first so the split "sticks" on both sides.
$clusterRng.Bold = $true
$clusterRng.Bold = $false
$perRng.Bold = $true
$perRng.Bold = $false

# --- 3) Drop the old "; clusters often extend into additional regions" ----
$deadRng = $d.Content
$deadRng.Find.Execute("; clusters often extend into additional regions", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$deadRng.Delete()

# --- 4) Re-seat the (hidden) "_GoBack" bookmark ----------------------------
# It used to span the whole sentence; now it should be an empty marker that
# sits right after "cluster" (and before the trailing ". ").
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$clusterRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $clusterRng)

Write-Output $d.Content.Text
